$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16 reuses the same formatting as row 15's A cell (bold, bordered, centered)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.193525908319875
$ws.Range("D16").Value = 0.8636495835432553
$ws.Range("E16").Value = 0.9527511789956422
$ws.Range("F16").Value = 1.193525908319875
$ws.Range("G16").Value = 0.8947636091639845
$ws.Range("H16").Value = 1.160794165004746
$ws.Range("I16").Value = 1.008482133597503
$ws.Range("J16").Value = 0.8636495835432553
$ws.Range("K16").Value = 0.9082003812694488
$ws.Range("L16").Value = 1.050863144794662
$ws.Range("M16").Value = 1.012327763104168
